$d = $word.ActiveDocument

$d.Content.Find.Execute("817×9=7353", $true, $false, $false, $false, $false, $true, 1, $false, "512×9=4608", 2)
$d.Content.Find.Execute("978×6=5868", $true, $false, $false, $false, $false, $true, 1, $false, "976×5=4880", 2)
$d.Content.Find.Execute("982×6=5892", $true, $false, $false, $false, $false, $true, 1, $false, "534×6=3204", 2)
$d.Content.Find.Execute("775×2=1550", $true, $false, $false, $false, $false, $true, 1, $false, "687×7=4809", 2)
$d.Content.Find.Execute("784×5=3920", $true, $false, $false, $false, $false, $true, 1, $false, "102×7=714", 2)
$d.Content.Find.Execute("392×2=784", $true, $false, $false, $false, $false, $true, 1, $false, "227×9=2043", 2)
$d.Content.Find.Execute("707×6=4242", $true, $false, $false, $false, $false, $true, 1, $false, "287×7=2009", 2)
$d.Content.Find.Execute("202×9=1818", $true, $false, $false, $false, $false, $true, 1, $false, "765×2=1530", 2)
$d.Content.Find.Execute("678×3=2034", $true, $false, $false, $false, $false, $true, 1, $false, "757×9=6813", 2)
$d.Content.Find.Execute("418×8=3344", $true, $false, $false, $false, $false, $true, 1, $false, "188×8=1504", 2)
$d.Content.Find.Execute("463×4=1852", $true, $false, $false, $false, $false, $true, 1, $false, "502×9=4518", 2)
$d.Content.Find.Execute("273×3=819", $true, $false, $false, $false, $false, $true, 1, $false, "477×4=1908", 2)
$d.Content.Find.Execute("838×8=6704", $true, $false, $false, $false, $false, $true, 1, $false, "366×5=1830", 2)
$d.Content.Find.Execute("194×9=1746", $true, $false, $false, $false, $false, $true, 1, $false, "850×2=1700", 2)
$d.Content.Find.Execute("502×8=4016", $true, $false, $false, $false, $false, $true, 1, $false, "156×5=780", 2)
$d.Content.Find.Execute("518×9=4662", $true, $false, $false, $false, $false, $true, 1, $false, "158×2=316", 2)
$d.Content.Find.Execute("102×8=816", $true, $false, $false, $false, $false, $true, 1, $false, "271×7=1897", 2)
$d.Content.Find.Execute("489×3=1467", $true, $false, $false, $false, $false, $true, 1, $false, "317×8=2536", 2)
$d.Content.Find.Execute("642×4=2568", $true, $false, $false, $false, $false, $true, 1, $false, "785×8=6280", 2)
$d.Content.Find.Execute("152×7=1064", $true, $false, $false, $false, $false, $true, 1, $false, "417×5=2085", 2)
$d.Content.Find.Execute("383×6=2298", $true, $false, $false, $false, $false, $true, 1, $false, "467×7=3269", 2)
$d.Content.Find.Execute("407×4=1628", $true, $false, $false, $false, $false, $true, 1, $false, "517×2=1034", 2)
$d.Content.Find.Execute("479×9=4311", $true, $false, $false, $false, $false, $true, 1, $false, "594×3=1782", 2)
$d.Content.Find.Execute("991×7=6937", $true, $false, $false, $false, $false, $true, 1, $false, "494×9=4446", 2)
$d.Content.Find.Execute("659×3=1977", $true, $false, $false, $false, $false, $true, 1, $false, "281×6=1686", 2)
